$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Update September_Details (R) and September_Date (S) columns for rows 45-174
# This reflects a new log entry inserted at row 45, shifting subsequent rows down by one.
$ws.Cells.Item(45, 18).Value = 'bal axis'
$ws.Cells.Item(45, 19).Value = '2024-09-21 07:56:12'
$ws.Cells.Item(46, 18).Value = 'bal axisbank axis'
$ws.Cells.Item(46, 19).Value = '2024-09-21 07:50:18'
$ws.Cells.Item(47, 18).Value = 'check the loan yo'
$ws.Cells.Item(47, 19).Value = '2024-09-20 15:37:11'
$ws.Cells.Item(48, 18).Value = 'balance your axis'
$ws.Cells.Item(48, 19).Value = '2024-09-20 08:05:28'
$ws.Cells.Item(49, 18).Value = 'bal axis'
$ws.Cells.Item(49, 19).Value = '2024-09-20 07:03:45'
$ws.Cells.Item(50, 18).Value = 'axis'
$ws.Cells.Item(50, 19).Value = '2024-09-20 06:57:43'
$ws.Cells.Item(51, 18).Value = 'dispute'
$ws.Cells.Item(51, 19).Value = '2024-09-19 22:46:00'
$ws.Cells.Item(52, 18).Value = 'tamilnadu disclose it anyone'
$ws.Cells.Item(52, 19).Value = '2024-09-19 22:41:11'
$ws.Cells.Item(53, 18).Value = 'dispute'
$ws.Cells.Item(53, 19).Value = '2024-09-19 22:33:39'
$ws.Cells.Item(54, 18).Value = 'dispute'
$ws.Cells.Item(54, 19).Value = '2024-09-19 22:27:16'
$ws.Cells.Item(55, 18).Value = 'your relationship'
$ws.Cells.Item(55, 19).Value = '2024-09-19 15:37:45'
$ws.Cells.Item(56, 18).Value = 'value discovery debit icici'
$ws.Cells.Item(56, 19).Value = '2024-09-19 14:34:40'
$ws.Cells.Item(57, 18).Value = 'debit'
$ws.Cells.Item(57, 19).Value = '2024-09-19 14:35:16'
$ws.Cells.Item(58, 18).Value = 'balance your axis'
$ws.Cells.Item(58, 19).Value = '2024-09-19 11:05:17'
$ws.Cells.Item(59, 18).Value = 'balance your axis'
$ws.Cells.Item(59, 19).Value = '2024-09-18 12:48:31'
$ws.Cells.Item(60, 18).Value = 'your relationship'
$ws.Cells.Item(60, 19).Value = '2024-09-18 10:29:06'
$ws.Cells.Item(61, 18).Value = 'balance your axis'
$ws.Cells.Item(61, 19).Value = '2024-09-18 10:28:28'
$ws.Cells.Item(62, 18).Value = 'axis'
$ws.Cells.Item(62, 19).Value = '2024-09-18 08:12:44'
$ws.Cells.Item(63, 18).Value = 'broker'
$ws.Cells.Item(63, 19).Value = '2024-09-18 04:09:58'
$ws.Cells.Item(64, 18).Value = 'balance your axis'
$ws.Cells.Item(64, 19).Value = '2024-09-17 13:07:16'
$ws.Cells.Item(65, 18).Value = 'dispute'
$ws.Cells.Item(65, 19).Value = '2024-09-16 12:53:44'
$ws.Cells.Item(66, 18).Value = 'money google icici'
$ws.Cells.Item(66, 19).Value = '2024-09-16 12:53:29'
$ws.Cells.Item(67, 18).Value = 'indusind'
$ws.Cells.Item(67, 19).Value = '2024-09-16 11:13:15'
$ws.Cells.Item(68, 18).Value = 'communication feedback'
$ws.Cells.Item(68, 19).Value = '2024-09-16 11:13:15'
$ws.Cells.Item(69, 18).Value = 'balance your axis'
$ws.Cells.Item(69, 19).Value = '2024-09-16 08:57:11'
$ws.Cells.Item(70, 18).Value = 'balance your axis'
$ws.Cells.Item(70, 19).Value = '2024-09-16 07:57:00'
$ws.Cells.Item(71, 18).Value = 'money google icici'
$ws.Cells.Item(71, 19).Value = '2024-09-15 21:06:00'
$ws.Cells.Item(72, 18).Value = 'adani icici'
$ws.Cells.Item(72, 19).Value = '2024-09-15 13:10:50'
$ws.Cells.Item(73, 18).Value = 'balance your axis'
$ws.Cells.Item(73, 19).Value = '2024-09-15 07:56:24'
$ws.Cells.Item(74, 18).Value = 'bal axisbank w axis'
$ws.Cells.Item(74, 19).Value = '2024-09-15 07:12:01'
$ws.Cells.Item(75, 18).Value = 'hdfc'
$ws.Cells.Item(75, 19).Value = '2024-09-14 21:25:23'
$ws.Cells.Item(76, 18).Value = 'change the'
$ws.Cells.Item(76, 19).Value = '2024-09-12 21:16:38'
$ws.Cells.Item(77, 18).Value = 'dispute'
$ws.Cells.Item(77, 19).Value = '2024-09-12 19:02:14'
$ws.Cells.Item(78, 18).Value = 'congrats limit icici'
$ws.Cells.Item(78, 19).Value = '2024-09-12 19:03:39'
$ws.Cells.Item(79, 18).Value = 'latest transaction pan'
$ws.Cells.Item(79, 19).Value = '2024-09-12 12:22:12'
$ws.Cells.Item(80, 18).Value = 'assistance'
$ws.Cells.Item(80, 19).Value = '2024-09-12 12:17:33'
$ws.Cells.Item(81, 18).Value = 'balance your axis'
$ws.Cells.Item(81, 19).Value = '2024-09-12 09:37:28'
$ws.Cells.Item(82, 18).Value = 'bal axisbank'
$ws.Cells.Item(82, 19).Value = '2024-09-12 00:54:39'
$ws.Cells.Item(83, 18).Value = 'your relationship'
$ws.Cells.Item(83, 19).Value = '2024-09-11 16:05:27'
$ws.Cells.Item(84, 18).Value = 'login internet personal share'
$ws.Cells.Item(84, 19).Value = '2024-09-11 14:16:45'
$ws.Cells.Item(85, 18).Value = 'balance your axis'
$ws.Cells.Item(85, 19).Value = '2024-09-11 12:45:33'
$ws.Cells.Item(86, 18).Value = 'balance your axis'
$ws.Cells.Item(86, 19).Value = '2024-09-11 09:45:01'
$ws.Cells.Item(87, 18).Value = 'axis'
$ws.Cells.Item(87, 19).Value = '2024-09-11 06:57:42'
$ws.Cells.Item(88, 18).Value = 'money google icici'
$ws.Cells.Item(88, 19).Value = '2024-09-10 20:42:12'
$ws.Cells.Item(89, 18).Value = 'dispute'
$ws.Cells.Item(89, 19).Value = '2024-09-10 20:42:34'
$ws.Cells.Item(90, 18).Value = 'reward points cash'
$ws.Cells.Item(90, 19).Value = '2024-09-10 19:43:35'
$ws.Cells.Item(91, 18).Value = 'balance your axis'
$ws.Cells.Item(91, 19).Value = '2024-09-10 13:32:42'
$ws.Cells.Item(92, 18).Value = 'ach indianesign bal axisbank'
$ws.Cells.Item(92, 19).Value = '2024-09-10 13:22:37'
$ws.Cells.Item(93, 18).Value = 'ach indianesign bal axisbank'
$ws.Cells.Item(93, 19).Value = '2024-09-10 13:22:37'
$ws.Cells.Item(94, 18).Value = 'balance your axis'
$ws.Cells.Item(94, 19).Value = '2024-09-10 11:21:40'
$ws.Cells.Item(95, 18).Value = 'your relationship'
$ws.Cells.Item(95, 19).Value = '2024-09-10 11:02:23'
$ws.Cells.Item(96, 18).Value = 'bank bal broker'
$ws.Cells.Item(96, 19).Value = '2024-09-09 19:59:02'
$ws.Cells.Item(97, 18).Value = 'beneficiary'
$ws.Cells.Item(97, 19).Value = '2024-09-09 15:48:10'
$ws.Cells.Item(98, 18).Value = 'beneficiary saravanan'
$ws.Cells.Item(98, 19).Value = '2024-09-09 14:52:20'
$ws.Cells.Item(99, 18).Value = 'bal axisbank'
$ws.Cells.Item(99, 19).Value = '2024-09-09 12:19:34'
$ws.Cells.Item(100, 18).Value = 'bal axisbank'
$ws.Cells.Item(100, 19).Value = '2024-09-09 12:19:33'
$ws.Cells.Item(101, 18).Value = 'dispute'
$ws.Cells.Item(101, 19).Value = '2024-09-09 12:17:30'
$ws.Cells.Item(102, 18).Value = 'bal axisbank'
$ws.Cells.Item(102, 19).Value = '2024-09-09 12:04:31'
$ws.Cells.Item(103, 18).Value = 'transfer freedom share anyone axis'
$ws.Cells.Item(103, 19).Value = '2024-09-09 11:56:19'
$ws.Cells.Item(104, 18).Value = 'corporate internet share'
$ws.Cells.Item(104, 19).Value = '2024-09-09 11:40:49'
$ws.Cells.Item(105, 18).Value = 'corporate internet share'
$ws.Cells.Item(105, 19).Value = '2024-09-09 11:39:30'
$ws.Cells.Item(106, 18).Value = 'bal axisbank'
$ws.Cells.Item(106, 19).Value = '2024-09-09 11:38:16'
$ws.Cells.Item(107, 18).Value = 'bal axisbank'
$ws.Cells.Item(107, 19).Value = '2024-09-09 11:38:16'
$ws.Cells.Item(108, 18).Value = 'bal axisbank'
$ws.Cells.Item(108, 19).Value = '2024-09-09 11:38:15'
$ws.Cells.Item(109, 18).Value = 'bal axisbank'
$ws.Cells.Item(109, 19).Value = '2024-09-09 11:38:15'
$ws.Cells.Item(110, 18).Value = 'corporate internet share'
$ws.Cells.Item(110, 19).Value = '2024-09-09 11:35:34'
$ws.Cells.Item(111, 18).Value = 'corporate internet share'
$ws.Cells.Item(111, 19).Value = '2024-09-09 11:32:23'
$ws.Cells.Item(112, 18).Value = 'ift anbu tpar'
$ws.Cells.Item(112, 19).Value = '2024-09-09 11:27:52'
$ws.Cells.Item(113, 18).Value = 'balance your axis'
$ws.Cells.Item(113, 19).Value = '2024-09-09 11:24:00'
$ws.Cells.Item(114, 18).Value = 'corporate internet share'
$ws.Cells.Item(114, 19).Value = '2024-09-09 11:21:43'
$ws.Cells.Item(115, 18).Value = 'corporate internet share'
$ws.Cells.Item(115, 19).Value = '2024-09-09 11:17:34'
$ws.Cells.Item(116, 18).Value = 'corporate internet share'
$ws.Cells.Item(116, 19).Value = '2024-09-09 11:15:51'
$ws.Cells.Item(117, 18).Value = 'corporate internet share'
$ws.Cells.Item(117, 19).Value = '2024-09-09 11:14:13'
$ws.Cells.Item(118, 18).Value = 'anbu tparty bal axisbank'
$ws.Cells.Item(118, 19).Value = '2024-09-09 11:13:37'
$ws.Cells.Item(119, 18).Value = 'corporate internet share'
$ws.Cells.Item(119, 19).Value = '2024-09-09 11:10:39'
$ws.Cells.Item(120, 18).Value = 'corporate internet share'
$ws.Cells.Item(120, 19).Value = '2024-09-09 11:07:31'
$ws.Cells.Item(121, 18).Value = 'corporate internet share'
$ws.Cells.Item(121, 19).Value = '2024-09-09 11:03:09'
$ws.Cells.Item(122, 18).Value = 'saravanan'
$ws.Cells.Item(122, 19).Value = '2024-09-09 10:43:11'
$ws.Cells.Item(123, 18).Value = 'balance your axis'
$ws.Cells.Item(123, 19).Value = '2024-09-09 08:10:16'
$ws.Cells.Item(124, 18).Value = 'ekalaivan'
$ws.Cells.Item(124, 19).Value = '2024-09-08 18:40:34'
$ws.Cells.Item(125, 18).Value = 'balance your axis'
$ws.Cells.Item(125, 19).Value = '2024-09-08 09:53:37'
$ws.Cells.Item(126, 18).Value = 'balance your axis'
$ws.Cells.Item(126, 19).Value = '2024-09-07 12:12:22'
$ws.Cells.Item(127, 18).Value = 'balance your axis'
$ws.Cells.Item(127, 19).Value = '2024-09-07 09:34:58'
$ws.Cells.Item(128, 18).Value = 'bal axis'
$ws.Cells.Item(128, 19).Value = '2024-09-07 08:46:40'
$ws.Cells.Item(129, 18).Value = 'axis'
$ws.Cells.Item(129, 19).Value = '2024-09-07 08:31:28'
$ws.Cells.Item(130, 18).Value = 'your relationship'
$ws.Cells.Item(130, 19).Value = '2024-09-06 12:23:25'
$ws.Cells.Item(131, 18).Value = 'balance your axis'
$ws.Cells.Item(131, 19).Value = '2024-09-06 09:55:31'
$ws.Cells.Item(132, 18).Value = 'beneficiary'
$ws.Cells.Item(132, 19).Value = '2024-09-05 17:13:56'
$ws.Cells.Item(133, 18).Value = 'coimbatore ramalinga'
$ws.Cells.Item(133, 19).Value = '2024-09-05 17:06:01'
$ws.Cells.Item(134, 18).Value = 'beneficiary'
$ws.Cells.Item(134, 19).Value = '2024-09-05 17:04:10'
$ws.Cells.Item(135, 18).Value = 'bal axisbank'
$ws.Cells.Item(135, 19).Value = '2024-09-05 16:52:25'
$ws.Cells.Item(136, 18).Value = 'share anyone axis'
$ws.Cells.Item(136, 19).Value = '2024-09-05 16:38:59'
$ws.Cells.Item(137, 18).Value = 'transfer anyone axis'
$ws.Cells.Item(137, 19).Value = '2024-09-05 16:35:58'
$ws.Cells.Item(138, 18).Value = 'share anyone axis'
$ws.Cells.Item(138, 19).Value = '2024-09-05 16:31:34'
$ws.Cells.Item(139, 18).Value = 'transfer'
$ws.Cells.Item(139, 19).Value = '2024-09-05 16:28:38'
$ws.Cells.Item(140, 18).Value = 'bal axisbank axis'
$ws.Cells.Item(140, 19).Value = '2024-09-05 16:26:56'
$ws.Cells.Item(141, 18).Value = 'bal axisbank'
$ws.Cells.Item(141, 19).Value = '2024-09-05 16:26:55'
$ws.Cells.Item(142, 18).Value = 'transfer'
$ws.Cells.Item(142, 19).Value = '2024-09-05 16:25:07'
$ws.Cells.Item(143, 18).Value = 'transfer'
$ws.Cells.Item(143, 19).Value = '2024-09-05 16:22:23'
$ws.Cells.Item(144, 18).Value = 'share anyone axis'
$ws.Cells.Item(144, 19).Value = '2024-09-05 16:06:05'
$ws.Cells.Item(145, 18).Value = 'internet bal axisbank'
$ws.Cells.Item(145, 19).Value = '2024-09-05 16:05:55'
$ws.Cells.Item(146, 18).Value = 'transfer share anyone axis'
$ws.Cells.Item(146, 19).Value = '2024-09-05 16:03:14'
$ws.Cells.Item(147, 18).Value = 'axis'
$ws.Cells.Item(147, 19).Value = '2024-09-05 15:57:15'
$ws.Cells.Item(148, 18).Value = 'your net internet'
$ws.Cells.Item(148, 19).Value = '2024-09-05 15:57:15'
$ws.Cells.Item(149, 18).Value = 'hear your feedback atm'
$ws.Cells.Item(149, 19).Value = '2024-09-05 14:21:08'
$ws.Cells.Item(150, 18).Value = 'axis bna'
$ws.Cells.Item(150, 19).Value = '2024-09-05 14:18:32'
$ws.Cells.Item(151, 18).Value = 'axis bna'
$ws.Cells.Item(151, 19).Value = '2024-09-05 14:13:16'
$ws.Cells.Item(152, 18).Value = 'axis bna'
$ws.Cells.Item(152, 19).Value = '2024-09-05 14:15:23'
$ws.Cells.Item(153, 18).Value = 'balance your axis'
$ws.Cells.Item(153, 19).Value = '2024-09-05 09:20:57'
$ws.Cells.Item(154, 18).Value = 'bal axis'
$ws.Cells.Item(154, 19).Value = '2024-09-05 09:06:25'
$ws.Cells.Item(155, 18).Value = 'broker'
$ws.Cells.Item(155, 19).Value = '2024-09-04 21:20:47'
$ws.Cells.Item(156, 18).Value = 'exclusive on axis'
$ws.Cells.Item(156, 19).Value = '2024-09-04 13:21:05'
$ws.Cells.Item(157, 18).Value = 'your corporate axis'
$ws.Cells.Item(157, 19).Value = '2024-09-04 11:46:10'
$ws.Cells.Item(158, 18).Value = 'balance your axis'
$ws.Cells.Item(158, 19).Value = '2024-09-04 08:14:16'
$ws.Cells.Item(159, 18).Value = 'axis'
$ws.Cells.Item(159, 19).Value = '2024-09-04 07:02:13'
$ws.Cells.Item(160, 18).Value = 'bal axisbank w axis'
$ws.Cells.Item(160, 19).Value = '2024-09-04 06:53:15'
$ws.Cells.Item(161, 18).Value = 'logging iob internet'
$ws.Cells.Item(161, 19).Value = '2024-09-03 20:09:12'
$ws.Cells.Item(162, 18).Value = 'password internet'
$ws.Cells.Item(162, 19).Value = '2024-09-03 20:05:31'
$ws.Cells.Item(163, 18).Value = 'logging iob internet'
$ws.Cells.Item(163, 19).Value = '2024-09-03 20:05:09'
$ws.Cells.Item(164, 18).Value = 'internet'
$ws.Cells.Item(164, 19).Value = '2024-09-03 19:58:18'
$ws.Cells.Item(165, 18).Value = 'login internet invalid'
$ws.Cells.Item(165, 19).Value = '2024-09-03 19:54:49'
$ws.Cells.Item(166, 18).Value = 'login internet invalid'
$ws.Cells.Item(166, 19).Value = '2024-09-03 19:56:17'
$ws.Cells.Item(167, 18).Value = 'corporate internet share'
$ws.Cells.Item(167, 19).Value = '2024-09-03 19:22:58'
$ws.Cells.Item(168, 18).Value = 'login sbi internet personal do not share anyone'
$ws.Cells.Item(168, 19).Value = '2024-09-03 19:17:10'
$ws.Cells.Item(169, 18).Value = 'login internet personal share'
$ws.Cells.Item(169, 19).Value = '2024-09-03 19:13:40'
$ws.Cells.Item(170, 18).Value = 'internet verify it'
$ws.Cells.Item(170, 19).Value = '2024-09-03 19:05:49'
$ws.Cells.Item(171, 18).Value = 'balance your axis'
$ws.Cells.Item(171, 19).Value = '2024-09-03 13:14:06'
$ws.Cells.Item(172, 18).Value = 'lounge'
$ws.Cells.Item(172, 19).Value = '2024-09-03 13:08:08'
$ws.Cells.Item(173, 18).Value = 'balance your axis'
$ws.Cells.Item(173, 19).Value = '2024-09-03 11:21:30'
$ws.Cells.Item(174, 18).Value = 'broker'
$ws.Cells.Item(174, 19).Value = '2024-09-01 22:35:38'

# Move the "Broadband" group label from row 182 to row 183
$ws.Cells.Item(182, 1).Value = $null
$ws.Cells.Item(183, 1).Value = "Broadband"
